$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BRTrains IDs")

# Insert two new rows after row 3 (pushes existing row 4+ down by 2)
$ws.Rows("4:5").Insert()

# Row 3 keeps only its ID value; clear the rest of the old data that has moved
$ws.Range("B3:E3").ClearContents()

# New row 4: BR Conflat A, now with an updated description and new ID 501
$ws.Range("A4").Value = 501
$ws.Range("B4").Value = "BR_Conflat_A"
$ws.Range("C4").Value = "BR Conflat A - Diagram 61/62"
$ws.Range("D4").Value = 1950
$ws.Range("E4").Value = "13 tons"

# New row 5: BR Conflat P, a brand new entry
$ws.Range("A5").Value = 502
$ws.Range("B5").Value = "BR_Conflat_P"
$ws.Range("C5").Value = "BR Conflat P - Diagram 60"
$ws.Range("D5").Value = 1959
$ws.Range("E5").Value = "12 tons"

$ws.Range("C6").Select()
